# Adding some more unit tests.
#
# This duplicates the "badPopValues_2" sheet into a new "badPopValues_4"
# sheet, fills in an "incomplete band" (E2/F2) on the new sheet, and
# leaves a review comment on that cell explaining the issue.

$wb = $excel.ActiveWorkbook

# --- Restore/record the cell selections that were left on the other
#     sheets before we add the new one (cosmetic, but cheap to match). ---
$wsNewPopValues   = $wb.Worksheets.Item("newPopValues")
$wsBadPopValues   = $wb.Worksheets.Item("badPopValues")
$wsBadPopValues2  = $wb.Worksheets.Item("badPopValues_2")
$wsBadPopValues3  = $wb.Worksheets.Item("badPopValues_3")

# --- Duplicate badPopValues_2 to create the new badPopValues_4 sheet. ---
$sourceSheet = $wsBadPopValues2
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sourceSheet.Copy($null, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "badPopValues_4"

# The copy also duplicated badPopValues_2's review comments; this new
# sheet should only carry the one new comment added below, so drop the
# comments that came along with the copy.
foreach ($addr in @("C1", "F7", "F21", "F25")) {
    $existing = $newSheet.Range($addr).Comment
    if ($existing -ne $null) {
        $existing.Delete()
    }
}

# Fill in the incomplete band that the comment below refers to.
$newSheet.Range("E2").Value = 2
$newSheet.Range("F2").Value = 98

# Leave a review comment flagging the incomplete band.
$newSheet.Range("E2").AddComment("Incomplete band - should go from 0-100")

# --- Cosmetic: leave the various sheets' selections where the author
#     left them. ---
$wsNewPopValues.Activate() | Out-Null
$wsNewPopValues.Range("K24").Select() | Out-Null

$wsBadPopValues.Activate() | Out-Null
$wsBadPopValues.Range("F7").Select() | Out-Null

$wsBadPopValues2.Activate() | Out-Null
$wsBadPopValues2.Range("L16").Select() | Out-Null

$wsBadPopValues3.Activate() | Out-Null
$wsBadPopValues3.Range("N17").Select() | Out-Null

$newSheet.Activate() | Out-Null
$newSheet.Range("O27").Select() | Out-Null
